$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $range = $d.Content
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1) Bold " que en el controlador dejé unas líneas comentadas" (leading space included),
#    leaving ", las cuales están relacionadas" un-bolded (this naturally splits the run).
$range = $d.Content
$found = $range.Find.Execute(" que en el controlador dejé unas líneas comentadas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $range.Font.Bold = 1
}

# 2) Fix "en el QUERY" -> "por los parámetros del QUERY" and rework the following sentence
ReplaceText "estuvieran dentro de la zona indicada en el QUERY)" "estuvieran dentro de la zona indicada por los parámetros del QUERY)"

ReplaceText "como no sé si esto se pudiera decidí dejar la opción que entre las 2 hacia menos iteraciones sin comentar y la que si usaba el cubo pero iteraba mucho más comentada" "como no sé si esto se permitiera, decidí dejar las 2 opciones, dejando sin comentar la opción que entre las 2 hacia menos iteraciones y la que si usaba el cubo pero iteraba mucho más quedó comentada"

# 3) Remove the _GoBack bookmark from its old location (it will be re-added later near "Preguntas")
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# 4) "la respuesta es una lista de los resultados" -> add "sencilla "
ReplaceText "la respuesta es una lista de los resultados" "la respuesta es una lista sencilla de los resultados"

# 5) "pero nunca se usa realmente después de ser creada." -> "pero nunca se usa realmente."
ReplaceText " pero nunca se usa realmente después de ser creada." " pero nunca se usa realmente."

# 6) Insert a new paragraph before "-El if que verifica..." (right after the $pushMessage paragraph)
$range = $d.Content
$found = $range.Find.Execute("-El ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertRange = $range.Paragraphs(1).Range
    $insertPoint = $d.Range($insertRange.Start, $insertRange.Start)
    $insertPoint.InsertBefore("-Comentarios que no ayudan a explicar el código, ni permitir intercambiar entre una manera de realizar la funcionalidad y otra (en caso de que se esté probando algo), ni a nada en realidad, simplemente están ahí y ya.`r")
    $insertPoint.ParagraphFormat.Alignment = 3
    $insertPoint.Font.Bold = 0
}

# 7) "en mi caso), no guardando result" -> "en mi caso), no se guarda la variable result"
ReplaceText " en mi caso), no guardando " " en mi caso), no se guarda la variable "

# 8) Insert "(no se crea realmente)" right after "result"
ReplaceText ", asignando el valor de la constante" " (no se crea realmente), asignando el valor de la constante"

# 9) "por fuera de la lógica, e hice" -> "por fuera de la lógica, ignoré los comentarios, e hice"
ReplaceText " por fuera de la lógica, " " por fuera de la lógica, ignoré los comentarios, "

# 10) "Preguntas" -> "Preguntas escritas" (bold run appended) + re-add _GoBack bookmark after it
$range = $d.Content
$found = $range.Find.Execute("Preguntas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($range.End, $range.End)
    $insertPoint.InsertAfter(" escritas")
    $insertPoint.Font.Bold = 1
    $insertPoint.Font.Size = 14
    $d.Bookmarks.Add("_GoBack", $insertPoint)
}
